# institutional_sites_metadata.xlsx - text/language updates to the
# "Description" column (B) on Sheet1, plus updated selection.
#
# The underlying shared-string table gets fully renumbered by the engine on
# save, so we don't try to match raw <si> indices from the diff - we just
# set the final cell text for every row whose wording actually changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: "15% Penetration Capacity (MW-AC)" column - expanded with an
# explanation of Rule 21 capacity rules (note the required trailing space).
$ws.Range("B15").Value = "15% Penetration Capacity (MW-AC) of the circuit serving the site. Importance dictated by Rule 21, the total combined nameplate capacity of all generation assets connected to each distribution circuit must not exceed 15% of the historical maximum load experienced on that circuit over the previous 18 months. "

# Row 19: Rooftop solar generation capacity -> nameplate solar generation
$ws.Range("B19").Value = "Rooftop solar nameplate solar generation potential (MW-AC)"

# Row 22: Low level parking lot estimate (MW-AC) - "generation capacity
# potential" -> "nameplate solar generation  potential" (double space kept
# exactly as authored)
$ws.Range("B22").Value = "Low level estimate of available parking lot area solar PV nameplate solar generation  potential (MW-AC)"

# Row 23: High level parking lot estimate (MW-AC), same wording change
$ws.Range("B23").Value = "High level estimate of available parking lot area solar PV nameplate solar generation  potential (MW-AC)"

# Row 24: Combined low estimate (MW-AC) - "generation capacity potential" ->
# "generation potential"
$ws.Range("B24").Value = "Combined total of the low estimate of parking lot generation potential and rooftop solar generation potential (MW-AC)"

# Row 25: Combined high estimate (MW-AC), same wording change
$ws.Range("B25").Value = "Combined total of the high estimate of parking lot generation potential and rooftop solar generation potential (MW-AC)"

# Row 32: Low level annual parking lot estimate (MWh AC) - drop "capacity"
$ws.Range("B32").Value = "Low level estimate of available parking lot area annual solar PV generation potential (MWh AC)"

# Row 33: High level annual parking lot estimate (MWh AC), same change
$ws.Range("B33").Value = "High level estimate of available parking lot area annual solar PV generation potential (MWh AC)"

# Row 35: Combined low annual estimate (MWh AC) -> "solar nameplate
# generation potential"
$ws.Range("B35").Value = "Combine total of the low estimate of parking lot and rooftop solar nameplate generation potential (MWh AC)"

# Row 36: Combined high annual estimate (MWh AC), same change
$ws.Range("B36").Value = "Combine total of the high estimate of parking lot and rooftop solar nameplate generation potential (MWh AC)"

# Update the sheet's active selection to match the new authored view state
# (was B21, now B25).
$ws.Range("B25").Select()
